$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 22.82060766666666
$ws.Cells.Item(2, 8).Value2 = 68.461823
$ws.Cells.Item(2, 9).Value2 = 0.10869226337084
$ws.Cells.Item(2, 10).Value2 = 0.1086922633708401
$ws.Cells.Item(2, 13).Value2 = 12.86269466666666
$ws.Cells.Item(2, 14).Value2 = 38.58808399999999
$ws.Cells.Item(2, 15).Value2 = 0.1337831063410017
$ws.Cells.Item(2, 16).Value2 = 0.1337831063410017
$ws.Cells.Item(2, 17).Value2 = 293.5345085241257
$ws.Cells.Item(2, 18).Value2 = 2641.810576717131
$ws.Cells.Item(2, 19).Value2 = 0.01454118862898526
$ws.Cells.Item(2, 20).Value2 = 0.01454118862898526

$ws.Cells.Item(3, 7).Value2 = 22.82060766666666
$ws.Cells.Item(3, 8).Value2 = 68.461823
$ws.Cells.Item(3, 9).Value2 = 0.10869226337084
$ws.Cells.Item(3, 10).Value2 = 0.1086922633708401
$ws.Cells.Item(3, 15).Value2 = 0.3593152390330854
$ws.Cells.Item(3, 16).Value2 = 0.3593152390330854
$ws.Cells.Item(3, 17).Value2 = 788.3762380727492
$ws.Cells.Item(3, 18).Value2 = 7095.386142654742
$ws.Cells.Item(3, 19).Value2 = 0.03905478659414047
$ws.Cells.Item(3, 20).Value2 = 0.03905478659414047

$ws.Cells.Item(4, 7).Value2 = 22.82060766666666
$ws.Cells.Item(4, 8).Value2 = 68.461823
$ws.Cells.Item(4, 9).Value2 = 0.10869226337084
$ws.Cells.Item(4, 10).Value2 = 0.1086922633708401
$ws.Cells.Item(4, 13).Value2 = 18.65324433333334
$ws.Cells.Item(4, 14).Value2 = 55.95973300000001
$ws.Cells.Item(4, 15).Value2 = 0.1940098117012772
$ws.Cells.Item(4, 16).Value2 = 0.1940098117012772
$ws.Cells.Item(4, 17).Value2 = 425.6783706414732
$ws.Cells.Item(4, 18).Value2 = 3831.105335773259
$ws.Cells.Item(4, 19).Value2 = 0.0210873655499623
$ws.Cells.Item(4, 20).Value2 = 0.0210873655499623

$ws.Cells.Item(5, 7).Value2 = 22.82060766666666
$ws.Cells.Item(5, 8).Value2 = 68.461823
$ws.Cells.Item(5, 9).Value2 = 0.10869226337084
$ws.Cells.Item(5, 10).Value2 = 0.1086922633708401
$ws.Cells.Item(5, 13).Value2 = 7.643308666666666
$ws.Cells.Item(5, 14).Value2 = 22.929926
$ws.Cells.Item(5, 15).Value2 = 0.07949699519803316
$ws.Cells.Item(5, 16).Value2 = 0.07949699519803316
$ws.Cells.Item(5, 17).Value2 = 174.4249483572331
$ws.Cells.Item(5, 18).Value2 = 1569.824535215098
$ws.Cells.Item(5, 19).Value2 = 0.008640708339255026
$ws.Cells.Item(5, 20).Value2 = 0.008640708339255027

$ws.Cells.Item(6, 7).Value2 = 22.82060766666666
$ws.Cells.Item(6, 8).Value2 = 68.461823
$ws.Cells.Item(6, 9).Value2 = 0.10869226337084
$ws.Cells.Item(6, 10).Value2 = 0.1086922633708401
$ws.Cells.Item(6, 13).Value2 = 22.43995333333334
$ws.Cells.Item(6, 14).Value2 = 67.31986000000001
$ws.Cells.Item(6, 15).Value2 = 0.2333948477266026
$ws.Cells.Item(6, 16).Value2 = 0.2333948477266026
$ws.Cells.Item(6, 17).Value2 = 512.0933710783089
$ws.Cells.Item(6, 18).Value2 = 4608.84033970478
$ws.Cells.Item(6, 19).Value2 = 0.02536821425849699
$ws.Cells.Item(6, 20).Value2 = 0.025368214258497

$ws.Cells.Item(7, 9).Value2 = 0.8059992924649123
$ws.Cells.Item(7, 10).Value2 = 0.8059992924649124
$ws.Cells.Item(7, 13).Value2 = 12.86269466666666
$ws.Cells.Item(7, 14).Value2 = 38.58808399999999
$ws.Cells.Item(7, 15).Value2 = 0.1337831063410017
$ws.Cells.Item(7, 16).Value2 = 0.1337831063410017
$ws.Cells.Item(7, 17).Value2 = 2176.683039318814
$ws.Cells.Item(7, 18).Value2 = 19590.14735386932
$ws.Cells.Item(7, 19).Value2 = 0.1078290890546055
$ws.Cells.Item(7, 20).Value2 = 0.1078290890546055

$ws.Cells.Item(8, 9).Value2 = 0.8059992924649123
$ws.Cells.Item(8, 10).Value2 = 0.8059992924649124
$ws.Cells.Item(8, 15).Value2 = 0.3593152390330854
$ws.Cells.Item(8, 16).Value2 = 0.3593152390330854
$ws.Cells.Item(8, 19).Value2 = 0.2896078284325277
$ws.Cells.Item(8, 20).Value2 = 0.2896078284325277

$ws.Cells.Item(9, 9).Value2 = 0.8059992924649123
$ws.Cells.Item(9, 10).Value2 = 0.8059992924649124
$ws.Cells.Item(9, 13).Value2 = 18.65324433333334
$ws.Cells.Item(9, 14).Value2 = 55.95973300000001
$ws.Cells.Item(9, 15).Value2 = 0.1940098117012772
$ws.Cells.Item(9, 16).Value2 = 0.1940098117012772
$ws.Cells.Item(9, 17).Value2 = 3156.585895944183
$ws.Cells.Item(9, 18).Value2 = 28409.27306349764
$ws.Cells.Item(9, 19).Value2 = 0.1563717709624803
$ws.Cells.Item(9, 20).Value2 = 0.1563717709624803

$ws.Cells.Item(10, 9).Value2 = 0.8059992924649123
$ws.Cells.Item(10, 10).Value2 = 0.8059992924649124
$ws.Cells.Item(10, 13).Value2 = 7.643308666666666
$ws.Cells.Item(10, 14).Value2 = 22.929926
$ws.Cells.Item(10, 15).Value2 = 0.07949699519803316
$ws.Cells.Item(10, 16).Value2 = 0.07949699519803316
$ws.Cells.Item(10, 17).Value2 = 1293.435067080177
$ws.Cells.Item(10, 18).Value2 = 11640.91560372159
$ws.Cells.Item(10, 19).Value2 = 0.06407452188270125
$ws.Cells.Item(10, 20).Value2 = 0.06407452188270127

$ws.Cells.Item(11, 9).Value2 = 0.8059992924649123
$ws.Cells.Item(11, 10).Value2 = 0.8059992924649124
$ws.Cells.Item(11, 13).Value2 = 22.43995333333334
$ws.Cells.Item(11, 14).Value2 = 67.31986000000001
$ws.Cells.Item(11, 15).Value2 = 0.2333948477266026
$ws.Cells.Item(11, 16).Value2 = 0.2333948477266026
$ws.Cells.Item(11, 17).Value2 = 3797.389823016792
$ws.Cells.Item(11, 18).Value2 = 34176.50840715112
$ws.Cells.Item(11, 19).Value2 = 0.1881160821325976
$ws.Cells.Item(11, 20).Value2 = 0.1881160821325977

$ws.Cells.Item(12, 7).Value2 = 4.587340666666667
$ws.Cells.Item(12, 8).Value2 = 13.762022
$ws.Cells.Item(12, 9).Value2 = 0.02184904307528146
$ws.Cells.Item(12, 10).Value2 = 0.02184904307528146
$ws.Cells.Item(12, 13).Value2 = 12.86269466666666
$ws.Cells.Item(12, 14).Value2 = 38.58808399999999
$ws.Cells.Item(12, 15).Value2 = 0.1337831063410017
$ws.Cells.Item(12, 16).Value2 = 0.1337831063410017
$ws.Cells.Item(12, 17).Value2 = 59.00556232731644
$ws.Cells.Item(12, 18).Value2 = 531.0500609458479
$ws.Cells.Item(12, 19).Value2 = 0.002923032853189507
$ws.Cells.Item(12, 20).Value2 = 0.002923032853189507

$ws.Cells.Item(13, 7).Value2 = 4.587340666666667
$ws.Cells.Item(13, 8).Value2 = 13.762022
$ws.Cells.Item(13, 9).Value2 = 0.02184904307528146
$ws.Cells.Item(13, 10).Value2 = 0.02184904307528146
$ws.Cells.Item(13, 15).Value2 = 0.3593152390330854
$ws.Cells.Item(13, 16).Value2 = 0.3593152390330854
$ws.Cells.Item(13, 17).Value2 = 158.4773915914336
$ws.Cells.Item(13, 18).Value2 = 1426.296524322902
$ws.Cells.Item(13, 19).Value2 = 0.007850694135238938
$ws.Cells.Item(13, 20).Value2 = 0.007850694135238938

$ws.Cells.Item(14, 7).Value2 = 4.587340666666667
$ws.Cells.Item(14, 8).Value2 = 13.762022
$ws.Cells.Item(14, 9).Value2 = 0.02184904307528146
$ws.Cells.Item(14, 10).Value2 = 0.02184904307528146
$ws.Cells.Item(14, 13).Value2 = 18.65324433333334
$ws.Cells.Item(14, 14).Value2 = 55.95973300000001
$ws.Cells.Item(14, 15).Value2 = 0.1940098117012772
$ws.Cells.Item(14, 16).Value2 = 0.1940098117012772
$ws.Cells.Item(14, 17).Value2 = 85.56878629556958
$ws.Cells.Item(14, 18).Value2 = 770.1190766601261
$ws.Cells.Item(14, 19).Value2 = 0.00423892873288845
$ws.Cells.Item(14, 20).Value2 = 0.00423892873288845

$ws.Cells.Item(15, 7).Value2 = 4.587340666666667
$ws.Cells.Item(15, 8).Value2 = 13.762022
$ws.Cells.Item(15, 9).Value2 = 0.02184904307528146
$ws.Cells.Item(15, 10).Value2 = 0.02184904307528146
$ws.Cells.Item(15, 13).Value2 = 7.643308666666666
$ws.Cells.Item(15, 14).Value2 = 22.929926
$ws.Cells.Item(15, 15).Value2 = 0.07949699519803316
$ws.Cells.Item(15, 16).Value2 = 0.07949699519803316
$ws.Cells.Item(15, 17).Value2 = 35.06246067448578
$ws.Cells.Item(15, 18).Value2 = 315.562146070372
$ws.Cells.Item(15, 19).Value2 = 0.00173693327243727
$ws.Cells.Item(15, 20).Value2 = 0.00173693327243727

$ws.Cells.Item(16, 7).Value2 = 4.587340666666667
$ws.Cells.Item(16, 8).Value2 = 13.762022
$ws.Cells.Item(16, 9).Value2 = 0.02184904307528146
$ws.Cells.Item(16, 10).Value2 = 0.02184904307528146
$ws.Cells.Item(16, 13).Value2 = 22.43995333333334
$ws.Cells.Item(16, 14).Value2 = 67.31986000000001
$ws.Cells.Item(16, 15).Value2 = 0.2333948477266026
$ws.Cells.Item(16, 16).Value2 = 0.2333948477266026
$ws.Cells.Item(16, 17).Value2 = 102.9397104841022
$ws.Cells.Item(16, 18).Value2 = 926.4573943569201
$ws.Cells.Item(16, 19).Value2 = 0.005099454081527297
$ws.Cells.Item(16, 20).Value2 = 0.005099454081527298

$ws.Cells.Item(17, 7).Value2 = 4.896139
$ws.Cells.Item(17, 8).Value2 = 14.688417
$ws.Cells.Item(17, 9).Value2 = 0.02331981853689061
$ws.Cells.Item(17, 10).Value2 = 0.02331981853689062
$ws.Cells.Item(17, 13).Value2 = 12.86269466666666
$ws.Cells.Item(17, 14).Value2 = 38.58808399999999
$ws.Cells.Item(17, 15).Value2 = 0.1337831063410017
$ws.Cells.Item(17, 16).Value2 = 0.1337831063410017
$ws.Cells.Item(17, 17).Value2 = 62.97754100255865
$ws.Cells.Item(17, 18).Value2 = 566.7978690230279
$ws.Cells.Item(17, 19).Value2 = 0.0031197977631737
$ws.Cells.Item(17, 20).Value2 = 0.003119797763173701

$ws.Cells.Item(18, 7).Value2 = 4.896139
$ws.Cells.Item(18, 8).Value2 = 14.688417
$ws.Cells.Item(18, 9).Value2 = 0.02331981853689061
$ws.Cells.Item(18, 10).Value2 = 0.02331981853689062
$ws.Cells.Item(18, 15).Value2 = 0.3593152390330854
$ws.Cells.Item(18, 16).Value2 = 0.3593152390330854
$ws.Cells.Item(18, 17).Value2 = 169.1453489005663
$ws.Cells.Item(18, 18).Value2 = 1522.308140105097
$ws.Cells.Item(18, 19).Value2 = 0.008379166171791028
$ws.Cells.Item(18, 20).Value2 = 0.00837916617179103

$ws.Cells.Item(19, 7).Value2 = 4.896139
$ws.Cells.Item(19, 8).Value2 = 14.688417
$ws.Cells.Item(19, 9).Value2 = 0.02331981853689061
$ws.Cells.Item(19, 10).Value2 = 0.02331981853689062
$ws.Cells.Item(19, 13).Value2 = 18.65324433333334
$ws.Cells.Item(19, 14).Value2 = 55.95973300000001
$ws.Cells.Item(19, 15).Value2 = 0.1940098117012772
$ws.Cells.Item(19, 16).Value2 = 0.1940098117012772
$ws.Cells.Item(19, 17).Value2 = 91.32887705696234
$ws.Cells.Item(19, 18).Value2 = 821.9598935126611
$ws.Cells.Item(19, 19).Value2 = 0.004524273603250101
$ws.Cells.Item(19, 20).Value2 = 0.004524273603250101

$ws.Cells.Item(20, 7).Value2 = 4.896139
$ws.Cells.Item(20, 8).Value2 = 14.688417
$ws.Cells.Item(20, 9).Value2 = 0.02331981853689061
$ws.Cells.Item(20, 10).Value2 = 0.02331981853689062
$ws.Cells.Item(20, 13).Value2 = 7.643308666666666
$ws.Cells.Item(20, 14).Value2 = 22.929926
$ws.Cells.Item(20, 15).Value2 = 0.07949699519803316
$ws.Cells.Item(20, 16).Value2 = 0.07949699519803316
$ws.Cells.Item(20, 17).Value2 = 37.42270165190466
$ws.Cells.Item(20, 18).Value2 = 336.804314867142
$ws.Cells.Item(20, 19).Value2 = 0.001853855502246198
$ws.Cells.Item(20, 20).Value2 = 0.001853855502246198

$ws.Cells.Item(21, 7).Value2 = 4.896139
$ws.Cells.Item(21, 8).Value2 = 14.688417
$ws.Cells.Item(21, 9).Value2 = 0.02331981853689061
$ws.Cells.Item(21, 10).Value2 = 0.02331981853689062
$ws.Cells.Item(21, 13).Value2 = 22.43995333333334
$ws.Cells.Item(21, 14).Value2 = 67.31986000000001
$ws.Cells.Item(21, 15).Value2 = 0.2333948477266026
$ws.Cells.Item(21, 16).Value2 = 0.2333948477266026
$ws.Cells.Item(21, 17).Value2 = 109.8691306735133
$ws.Cells.Item(21, 18).Value2 = 988.82217606162
$ws.Cells.Item(21, 19).Value2 = 0.005442725496429589
$ws.Cells.Item(21, 20).Value2 = 0.00544272549642959

$ws.Cells.Item(22, 7).Value2 = 8.427551666666668
$ws.Cells.Item(22, 8).Value2 = 25.282655
$ws.Cells.Item(22, 9).Value2 = 0.04013958255207557
$ws.Cells.Item(22, 10).Value2 = 0.04013958255207558
$ws.Cells.Item(22, 13).Value2 = 12.86269466666666
$ws.Cells.Item(22, 14).Value2 = 38.58808399999999
$ws.Cells.Item(22, 15).Value2 = 0.1337831063410017
$ws.Cells.Item(22, 16).Value2 = 0.1337831063410017
$ws.Cells.Item(22, 17).Value2 = 108.4010238758911
$ws.Cells.Item(22, 18).Value2 = 975.6092148830199
$ws.Cells.Item(22, 19).Value2 = 0.005369998041047743
$ws.Cells.Item(22, 20).Value2 = 0.005369998041047744

$ws.Cells.Item(23, 7).Value2 = 8.427551666666668
$ws.Cells.Item(23, 8).Value2 = 25.282655
$ws.Cells.Item(23, 9).Value2 = 0.04013958255207557
$ws.Cells.Item(23, 10).Value2 = 0.04013958255207558
$ws.Cells.Item(23, 15).Value2 = 0.3593152390330854
$ws.Cells.Item(23, 16).Value2 = 0.3593152390330854
$ws.Cells.Item(23, 17).Value2 = 291.143933420984
$ws.Cells.Item(23, 18).Value2 = 2620.295400788855
$ws.Cells.Item(23, 19).Value2 = 0.0144227636993873
$ws.Cells.Item(23, 20).Value2 = 0.0144227636993873

$ws.Cells.Item(24, 7).Value2 = 8.427551666666668
$ws.Cells.Item(24, 8).Value2 = 25.282655
$ws.Cells.Item(24, 9).Value2 = 0.04013958255207557
$ws.Cells.Item(24, 10).Value2 = 0.04013958255207558
$ws.Cells.Item(24, 13).Value2 = 18.65324433333334
$ws.Cells.Item(24, 14).Value2 = 55.95973300000001
$ws.Cells.Item(24, 15).Value2 = 0.1940098117012772
$ws.Cells.Item(24, 16).Value2 = 0.1940098117012772
$ws.Cells.Item(24, 17).Value2 = 157.201180370124
$ws.Cells.Item(24, 18).Value2 = 1414.810623331115
$ws.Cells.Item(24, 19).Value2 = 0.007787472852696052
$ws.Cells.Item(24, 20).Value2 = 0.007787472852696053

$ws.Cells.Item(25, 7).Value2 = 8.427551666666668
$ws.Cells.Item(25, 8).Value2 = 25.282655
$ws.Cells.Item(25, 9).Value2 = 0.04013958255207557
$ws.Cells.Item(25, 10).Value2 = 0.04013958255207558
$ws.Cells.Item(25, 13).Value2 = 7.643308666666666
$ws.Cells.Item(25, 14).Value2 = 22.929926
$ws.Cells.Item(25, 15).Value2 = 0.07949699519803316
$ws.Cells.Item(25, 16).Value2 = 0.07949699519803316
$ws.Cells.Item(25, 17).Value2 = 64.41437869261445
$ws.Cells.Item(25, 18).Value2 = 579.72940823353
$ws.Cells.Item(25, 19).Value2 = 0.003190976201393407
$ws.Cells.Item(25, 20).Value2 = 0.003190976201393408

$ws.Cells.Item(26, 7).Value2 = 8.427551666666668
$ws.Cells.Item(26, 8).Value2 = 25.282655
$ws.Cells.Item(26, 9).Value2 = 0.04013958255207557
$ws.Cells.Item(26, 10).Value2 = 0.04013958255207558
$ws.Cells.Item(26, 13).Value2 = 22.43995333333334
$ws.Cells.Item(26, 14).Value2 = 67.31986000000001
$ws.Cells.Item(26, 15).Value2 = 0.2333948477266026
$ws.Cells.Item(26, 16).Value2 = 0.2333948477266026
$ws.Cells.Item(26, 17).Value2 = 189.1138661142556
$ws.Cells.Item(26, 18).Value2 = 1702.0247950283
$ws.Cells.Item(26, 19).Value2 = 0.009368371757551073
$ws.Cells.Item(26, 20).Value2 = 0.009368371757551074
